$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the SPI SDI/SDO pin assignments between row 17 (RC6) and row 18 (RC5)
$i17 = $ws.Range("I17").Value()
$j17 = $ws.Range("J17").Value()
$k17 = $ws.Range("K17").Value()

$i18 = $ws.Range("I18").Value()
$j18 = $ws.Range("J18").Value()
$k18 = $ws.Range("K18").Value()

$ws.Range("I17").Value = $i18
$ws.Range("J17").Value = $j18
$ws.Range("K17").Value = $k18

$ws.Range("I18").Value = $i17
$ws.Range("J18").Value = $j17
$ws.Range("K18").Value = $k17

# Update the active selection to K19
$ws.Range("K19").Select()
